$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark (paragraph ending "Thanh toan: ") ---
$d.Bookmarks.Item("_GoBack").Delete()

# --- 2. Append new paragraphs at the end of the document body ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)

# two new blank paragraphs (move past each inserted paragraph mark so the
# range tracks forward into the freshly-created empty paragraph)
$r.InsertParagraphAfter()
$r.Move(1, 1) | Out-Null
$r.InsertParagraphAfter()
$r.Move(1, 1) | Out-Null

# third new paragraph -- this is the one that receives the long text
$r.InsertParagraphAfter()
$r.Move(1, 1) | Out-Null

$r.InsertAfter('Merchant_site_code + '' ')
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
$r.Collapse(0)

$r.InsertAfter(''' + return_url + '' '' + receiver + '' '' + transaction_info + '' '' + order_code + '' '' + price + '' '' + currency + '' '' + quantity + '' '' + tax + '' '' + discount + '' '' + fee_cal + '' '' + fee_shipping + '' '' + Skype: hotrokythuat02@nganluong.vn Mail: hotrokythuat@nganluong.vn NganLuong.vn – Tích hợp tiêu chuẩn – Checkout version 2.0 order_description + '' '' + buyer_info + '' '' + affiliate_code + '' '' + secure_pass '')`')
$r.Collapse(0)

Write-Host "done"
